$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$wsALC.Range("H2").Value = 231
$wsALC.Range("I2").Value = 266.5
$wsALC.Range("J2").Value = 160
$wsALC.Range("K2").Value = 266.5
$wsALC.Range("L2").Value = 160
$wsALC.Range("M2").Value = -153.5
$wsALC.Range("N2").Value = -386
$wsALC.Range("H53").Value = 298.73914
$wsALC.Range("I53").Value = 134.83333
$wsALC.Range("J53").Value = 356.58823
$wsALC.Range("K53").Value = 134.83333
$wsALC.Range("L53").Value = 356.58823
$wsALC.Range("M53").Value = 502.16667
$wsALC.Range("N53").Value = -1630.58823
$wsALC.Range("H86").Value = 5072.5835
$wsALC.Range("I86").Value = 3840.2
$wsALC.Range("J86").Value = 5952.857
$wsALC.Range("K86").Value = 3840.2
$wsALC.Range("L86").Value = 5952.857
$wsALC.Range("M86").Value = -2717.2
$wsALC.Range("N86").Value = -8198.857
$wsALC.Range("H89").Value = 5072.5835
$wsALC.Range("I89").Value = 3840.2
$wsALC.Range("J89").Value = 5952.857
$wsALC.Range("K89").Value = 19201
$wsALC.Range("L89").Value = 29764.285
$wsALC.Range("M89").Value = -13585
$wsALC.Range("N89").Value = -40996.285
$wsALC.Range("H111").Value = 3464.923
$wsALC.Range("I111").Value = 3218.889
$wsALC.Range("K111").Value = 9656.667000000001
$wsALC.Range("M111").Value = -6589.667000000001
$wsALC.Range("H130").Value = 44997.5
$wsALC.Range("J130").Value = 44997.5
$wsALC.Range("L130").Value = 44997.5
$wsALC.Range("N130").Value = -55037.5
$wsALC.Range("H134").Value = 39615.383
$wsALC.Range("J134").Value = 39615.383
$wsALC.Range("L134").Value = 39615.383
$wsALC.Range("N134").Value = -49755.383
$wsALC.Range("H136").Value = 39666.668
$wsALC.Range("J136").Value = 39666.668
$wsALC.Range("L136").Value = 39666.668
$wsALC.Range("N136").Value = -49866.668
$wsALC.Range("H137").Value = 5066.6875
$wsALC.Range("I137").Value = 5078.6665
$wsALC.Range("J137").Value = 5030.75
$wsALC.Range("K137").Value = 15235.9995
$wsALC.Range("L137").Value = 15092.25
$wsALC.Range("M137").Value = -12685.9995
$wsALC.Range("N137").Value = -20192.25
$wsALC.Range("H138").Value = 3039.795
$wsALC.Range("I138").Value = 1503
$wsALC.Range("K138").Value = 4509
$wsALC.Range("M138").Value = 631
$wsALC.Range("H139").Value = 49597.5
$wsALC.Range("J139").Value = 49597.5
$wsALC.Range("L139").Value = 49597.5
$wsALC.Range("N139").Value = -59877.5

# --- ARM ---
$wsARM.Range("H74").Value = 1758.5333
$wsARM.Range("I74").Value = 1741.2858
$wsARM.Range("K74").Value = 1741.2858
$wsARM.Range("M74").Value = -867.2858000000001
$wsARM.Range("H77").Value = 1758.5333
$wsARM.Range("I77").Value = 1741.2858
$wsARM.Range("K77").Value = 8706.429
$wsARM.Range("M77").Value = -4338.429

# --- BSM ---
$wsBSM.Range("H134").Value = 1794.5428
$wsBSM.Range("I134").Value = 1198.8948
$wsBSM.Range("K134").Value = 3596.6844
$wsBSM.Range("M134").Value = -1061.6844

# --- CRP ---
$wsCRP.Range("H22").Value = 849.1667
$wsCRP.Range("I22").Value = 732.6667
$wsCRP.Range("K22").Value = 732.6667
$wsCRP.Range("M22").Value = -382.6667
$wsCRP.Range("H132").Value = 8399.200000000001
$wsCRP.Range("I132").Value = 0
$wsCRP.Range("K132").Value = 0
$wsCRP.Range("M132").ClearContents()
$wsCRP.Range("H134").Value = 2081.889
$wsCRP.Range("I134").Value = 2089.8
$wsCRP.Range("J134").Value = 2042.3334
$wsCRP.Range("K134").Value = 6269.400000000001
$wsCRP.Range("L134").Value = 6127.0002
$wsCRP.Range("M134").Value = -3734.400000000001
$wsCRP.Range("N134").Value = -11197.0002

# --- CUL ---
$wsCUL.Range("H22").Value = 66666916
$wsCUL.Range("I22").Value = 344
$wsCUL.Range("J22").Value = 166666770
$wsCUL.Range("K22").Value = 1032
$wsCUL.Range("L22").Value = 500000310
$wsCUL.Range("M22").Value = -863
$wsCUL.Range("N22").Value = -500000648
$wsCUL.Range("H27").Value = 66666916
$wsCUL.Range("I27").Value = 344
$wsCUL.Range("J27").Value = 166666770
$wsCUL.Range("K27").Value = 1032
$wsCUL.Range("L27").Value = 500000310
$wsCUL.Range("M27").Value = -930
$wsCUL.Range("N27").Value = -500000514
$wsCUL.Range("H93").Value = 7499.25
$wsCUL.Range("I93").Value = 0
$wsCUL.Range("J93").Value = 7499.25
$wsCUL.Range("K93").Value = 0
$wsCUL.Range("L93").Value = 22497.75
$wsCUL.Range("M93").ClearContents()
$wsCUL.Range("N93").Value = -26241.75
$wsCUL.Range("H107").Value = 2690.4075
$wsCUL.Range("J107").Value = 2685.4707
$wsCUL.Range("L107").Value = 8056.4121
$wsCUL.Range("N107").Value = -11896.4121
$wsCUL.Range("H113").Value = 1483.875
$wsCUL.Range("I113").Value = 1660.25
$wsCUL.Range("J113").Value = 1307.5
$wsCUL.Range("K113").Value = 4980.75
$wsCUL.Range("L113").Value = 3922.5
$wsCUL.Range("M113").Value = -2810.75
$wsCUL.Range("N113").Value = -8262.5
$wsCUL.Range("H122").Value = 1229.3846
$wsCUL.Range("I122").Value = 609.6
$wsCUL.Range("K122").Value = 5486.400000000001
$wsCUL.Range("M122").Value = -3036.400000000001

# --- GSM ---
$wsGSM.Range("H132").Value = 5386.0713
$wsGSM.Range("I132").Value = 4645.851
$wsGSM.Range("K132").Value = 13937.553
$wsGSM.Range("M132").Value = -11407.553

# --- LTW ---
$wsLTW.Range("H4").Value = 0
$wsLTW.Range("I4").Value = 0
$wsLTW.Range("K4").Value = 0
$wsLTW.Range("M4").ClearContents()
$wsLTW.Range("H22").Value = 15152554
$wsLTW.Range("J22").Value = 1994.5
$wsLTW.Range("L22").Value = 1994.5
$wsLTW.Range("N22").Value = -2584.5
$wsLTW.Range("H27").Value = 15152554
$wsLTW.Range("J27").Value = 1994.5
$wsLTW.Range("L27").Value = 1994.5
$wsLTW.Range("N27").Value = -2208.5
$wsLTW.Range("H28").Value = 0
$wsLTW.Range("I28").Value = 0
$wsLTW.Range("K28").Value = 0
$wsLTW.Range("M28").ClearContents()
$wsLTW.Range("H32").Value = 2422.6
$wsLTW.Range("I32").Value = 2422.6
$wsLTW.Range("K32").Value = 2422.6
$wsLTW.Range("M32").Value = -2105.6
$wsLTW.Range("H33").Value = 9005
$wsLTW.Range("J33").Value = 10000
$wsLTW.Range("L33").Value = 10000
$wsLTW.Range("N33").Value = -10580
$wsLTW.Range("H37").Value = 0
$wsLTW.Range("I37").Value = 0
$wsLTW.Range("K37").Value = 0
$wsLTW.Range("M37").ClearContents()
$wsLTW.Range("H123").Value = 29166.666
$wsLTW.Range("J123").Value = 29166.666
$wsLTW.Range("L123").Value = 29166.666
$wsLTW.Range("N123").Value = -38966.666
$wsLTW.Range("H136").Value = 4054
$wsLTW.Range("I136").Value = 3552.5
$wsLTW.Range("J136").Value = 4627.143
$wsLTW.Range("K136").Value = 10657.5
$wsLTW.Range("L136").Value = 13881.429
$wsLTW.Range("M136").Value = -8107.5
$wsLTW.Range("N136").Value = -18981.429
$wsLTW.Range("H140").Value = 63398.6
$wsLTW.Range("J140").Value = 74248.25
$wsLTW.Range("L140").Value = 74248.25
$wsLTW.Range("N140").Value = -84608.25

# --- WVR ---
$wsWVR.Range("H40").Value = 0
$wsWVR.Range("J40").Value = 0
$wsWVR.Range("L40").Value = 0
$wsWVR.Range("N40").ClearContents()
$wsWVR.Range("H132").Value = 4064.0625
$wsWVR.Range("I132").Value = 3693.0435
$wsWVR.Range("K132").Value = 11079.1305
$wsWVR.Range("M132").Value = -8549.130500000001

Write-Output "applied all changes"